$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.7059925093632958
$ws1.Range("C2").Value = 0.6656626506024096
$ws1.Range("D2").Value = 0.8277153558052435
$ws1.Range("E2").Value = 0.7378964941569283
$ws1.Range("F2").Value = 0.7892857142857143
$ws1.Range("G2").Value = 0.8200371057513914
$ws1.Range("H2").Value = 0.818462525775365
$ws1.Range("I2").Value = 442
$ws1.Range("J2").Value = 222
$ws1.Range("K2").Value = 312
$ws1.Range("L2").Value = 92

# --- Classification Report sheet ---
$ws2 = $wb.Worksheets.Item("Classification Report")
$ws2.Range("B2").Value = 0.7722772277227723
$ws2.Range("C2").Value = 0.5842696629213483
$ws2.Range("D2").Value = 0.6652452025586354

$ws2.Range("B3").Value = 0.6656626506024096
$ws2.Range("C3").Value = 0.8277153558052435
$ws2.Range("D3").Value = 0.7378964941569283

$ws2.Range("B4").Value = 0.7059925093632958
$ws2.Range("C4").Value = 0.7059925093632958
$ws2.Range("D4").Value = 0.7059925093632958
$ws2.Range("E4").Value = 0.7059925093632958

$ws2.Range("B5").Value = 0.718969939162591
$ws2.Range("C5").Value = 0.7059925093632959
$ws2.Range("D5").Value = 0.7015708483577818

$ws2.Range("B6").Value = 0.7189699391625909
$ws2.Range("C6").Value = 0.7059925093632958
$ws2.Range("D6").Value = 0.7015708483577818

# --- Confusion Matrix sheet ---
$ws3 = $wb.Worksheets.Item("Confusion Matrix")
$ws3.Range("B2").Value = 312
$ws3.Range("C2").Value = 222
$ws3.Range("B3").Value = 92
$ws3.Range("C3").Value = 442
